# Edit script: applies the two changes captured in the target diff.
#
# 1) Slide 16's table (shape 3, a graphicFrame) switches to a different
#    table style (tableStyleId GUID change).
# 2) The deck's primary theme (ppt/theme/theme1.xml, applied via the
#    slide master) has its 12 theme colors swapped from the custom
#    "Integral" palette to the stock Office default palette (and picks
#    up the "Office Theme" identity that used to live in theme2.xml).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{646F234E-5681-4EEE-AB00-5A31D2AC679D}")

# --- 2. Recolor the theme (Integral -> stock Office palette) --------------
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = "0x000000"   # dk1
$cs.Colors(2).RGB  = "0xFFFFFF"   # lt1
$cs.Colors(3).RGB  = "0x44546A"   # dk2
$cs.Colors(4).RGB  = "0xE7E6E6"   # lt2
$cs.Colors(5).RGB  = "0x5B9BD5"   # accent1
$cs.Colors(6).RGB  = "0xED7D31"   # accent2
$cs.Colors(7).RGB  = "0xA5A5A5"   # accent3
$cs.Colors(8).RGB  = "0xFFC000"   # accent4
$cs.Colors(9).RGB  = "0x4472C4"   # accent5
$cs.Colors(10).RGB = "0x70AD47"   # accent6
$cs.Colors(11).RGB = "0x0563C1"   # hyperlink
$cs.Colors(12).RGB = "0x954F72"   # followed hyperlink
